# Design basis sheets refactored
#
# The workbook's three equipment sheets are renamed with a "VTUS88BP"
# project/tag suffix, and the "MCC CUM PLC" sheet (now "MCC CUM PLC VTUS88BP")
# becomes the active/selected sheet instead of "COVER".

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("MCC").Name = "MCC VTUS88BP"
$wb.Worksheets.Item("PCC").Name = "PCC VTUS88BP"
$wb.Worksheets.Item("MCC CUM PLC").Name = "MCC CUM PLC VTUS88BP"

# Make "MCC CUM PLC VTUS88BP" the active tab (was "COVER").
$wb.Worksheets.Item("MCC CUM PLC VTUS88BP").Activate()
